$d = $word.ActiveDocument

# Table 1 value refresh: replace each old statistic with its updated value.
# Each (old, new) pair below is unique across the document, so a plain
# Find/Replace (whole document, match case, no wildcards) is unambiguous.
$replacements = @(
    @('All (n=321073)', 'All (n=274989)'),
    @('Private Insurance Group (n=189036)', 'Private Insurance Group (n=158402)'),
    @('Self-pay Group (n=43373)', 'Self-pay Group (n=39468)'),
    @('41.69±13.33', '40.34±13.35'),
    @('43.09±13.18', '41.66±13.31'),
    @('37.14±12.34', '36.33±12.23'),
    @('9942 (3.10)', '9480 (3.45)'),
    @('6207 (3.28)', '5870 (3.71)'),
    @('986 (2.27)', '960 (2.43)'),
    @('27678 (8.62)', '23574 (8.57)'),
    @('12024 (6.36)', '9980 (6.30)'),
    @('4642 (10.70)', '4163 (10.55)'),
    @('61841 (19.26)', '55990 (20.36)'),
    @('22704 (12.01)', '19844 (12.53)'),
    @('14162 (32.65)', '13210 (33.47)'),
    @('1820 (0.57)', '1575 (0.57)'),
    @('811 (0.43)', '684 (0.43)'),
    @('224 (0.52)', '207 (0.52)'),
    @('12486 (3.89)', '11339 (4.12)'),
    @('6213 (3.29)', '5519 (3.48)'),
    @('2238 (5.16)', '2095 (5.31)'),
    @('17234 (5.37)', '15995 (5.82)'),
    @('11147 (5.90)', '10243 (6.47)'),
    @('1948 (4.49)', '1887 (4.78)'),
    @('190072 (59.20)', '157036 (57.11)'),
    @('129930 (68.73)', '106262 (67.08)'),
    @('19173 (44.20)', '16946 (42.94)'),
    @('150169 (46.77)', '131408 (47.79)'),
    @('89295 (47.24)', '76851 (48.52)'),
    @('16445 (37.92)', '15205 (38.52)'),
    @('169783 (52.88)', '142467 (51.81)'),
    @('99066 (52.41)', '80882 (51.06)'),
    @('26755 (61.69)', '24090 (61.04)'),
    @('1121 (0.35)', '1114 (0.41)'),
    @('675 (0.36)', '669 (0.42)'),
    @('82920 (25.83)', '71903 (26.15)'),
    @('35451 (18.75)', '29745 (18.78)'),
    @('16114 (37.15)', '14675 (37.18)'),
    @('79085 (24.63)', '67128 (24.41)'),
    @('42906 (22.70)', '35407 (22.35)'),
    @('11997 (27.66)', '10832 (27.45)'),
    @('82093 (25.57)', '69846 (25.40)'),
    @('52301 (27.67)', '43568 (27.50)'),
    @('9607 (22.15)', '8772 (22.23)'),
    @('76975 (23.97)', '66112 (24.04)'),
    @('58378 (30.88)', '49682 (31.36)'),
    @('5655 (13.04)', '5189 (13.15)'),
    @('32121 (10.00)', '27877 (10.14)'),
    @('17903 (9.47)', '15243 (9.62)'),
    @('5106 (11.77)', '4675 (11.85)'),
    @('114867 (35.78)', '103415 (37.61)'),
    @('69436 (36.73)', '61848 (39.04)'),
    @('16020 (36.94)', '15002 (38.01)'),
    @('174085 (54.22)', '143697 (52.26)'),
    @('101697 (53.80)', '81311 (51.33)'),
    @('22247 (51.29)', '19791 (50.14)'),
    @('58818 (18.32)', '48168 (17.52)'),
    @('38226 (20.22)', '30836 (19.47)'),
    @('6385 (14.72)', '5735 (14.53)'),
    @('66027 (20.56)', '56761 (20.64)'),
    @('40759 (21.56)', '34393 (21.71)'),
    @('5572 (12.85)', '5217 (13.22)'),
    @('115219 (35.89)', '97440 (35.43)'),
    @('65081 (34.43)', '53443 (33.74)'),
    @('23608 (54.43)', '21065 (53.37)'),
    @('81009 (25.23)', '72620 (26.41)'),
    @('44970 (23.79)', '39730 (25.08)'),
    @('7808 (18.00)', '7451 (18.88)'),
    @('229819 (71.58)', '229819 (83.57)'),
    @('136044 (71.97)', '136044 (85.89)'),
    @('33575 (77.41)', '33575 (85.07)'),
    @('40947 (12.75)', '40947 (14.89)'),
    @('19983 (10.57)', '19983 (12.62)'),
    @('5451 (12.57)', '5451 (13.81)'),
    @('50307 (15.67)', '4223 (1.54)'),
    @('33009 (17.46)', '2375 (1.50)'),
    @('4347 (10.02)', '442 (1.12)'),
    @('174816 (54.45)', '150145 (54.60)'),
    @('104335 (55.19)', '87251 (55.08)'),
    @('25525 (58.85)', '23273 (58.97)'),
    @('123812 (38.56)', '107709 (39.17)'),
    @('73483 (38.87)', '63047 (39.80)'),
    @('15609 (35.99)', '14359 (36.38)'),
    @('19338 (6.02)', '14885 (5.41)'),
    @('9780 (5.17)', '7110 (4.49)'),
    @('1912 (4.41)', '1580 (4.00)'),
    @('3107 (0.97)', '2250 (0.82)'),
    @('1438 (0.76)', '994 (0.63)'),
    @('327 (0.75)', '256 (0.65)'),
    @('290181 (90.38)', '251390 (91.42)'),
    @('173376 (91.72)', '147074 (92.85)'),
    @('40465 (93.30)', '37103 (94.01)'),
    @('23045 (7.18)', '17539 (6.38)'),
    @('12076 (6.39)', '8687 (5.48)'),
    @('2167 (5.00)', '1759 (4.46)'),
    @('6239 (1.94)', '4865 (1.77)'),
    @('2843 (1.50)', '2114 (1.33)'),
    @('568 (1.31)', '471 (1.19)'),
    @('1608 (0.50)', '1195 (0.43)'),
    @('741 (0.39)', '527 (0.33)'),
    @('0.37±1.20', '0.33±1.15'),
    @('0.31±1.10', '0.27±1.06'),
    @('0.19±0.71', '0.18±0.67'),
    @('289 (0.09)', '213 (0.08)'),
    @('98 (0.05)', '67 (0.04)'),
    @('28 (0.06)', '21 (0.05)'),
    @('83864 (26.12)', '68404 (24.88)'),
    @('48383 (25.59)', '37805 (23.87)'),
    @('10085 (23.25)', '9080 (23.01)'),
    @('29708 (9.25)', '25110 (9.13)'),
    @('15695 (8.30)', '13022 (8.22)'),
    @('4110 (9.48)', '3654 (9.26)'),
    @('1269 (0.40)', '1187 (0.43)'),
    @('705 (0.37)', '656 (0.41)'),
    @('191 (0.44)', '184 (0.47)'),
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null
}

# The two cells below both originally read "173 (0.40)" so a global text
# Find/Replace would be ambiguous; address them positionally instead.
$tbl = $d.Tables(1)
$tbl.Cell(14, 4).Range.Text = "173 (0.44)"
$tbl.Cell(42, 4).Range.Text = "135 (0.34)"

Write-Host "Applied" $replacements.Count "Find/Replace updates plus 2 positional cell updates."